$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Evaluation" to "Sheet1"
$ws.Name = "Sheet1"

# Update Cohere_Translation (column C) text and BLEU_Score (column D) values
$ws.Cells.Item(8, 3).Value = "আর্জেন্টিনা জানা মানে সে বেস্ট পোলো টিম এবং খেলোয়াড়দের থাকা, বিশ্বের সেরাদের মধ্যে একটি।"
$ws.Cells.Item(8, 4).Value = 0.02486954214312595

$ws.Cells.Item(12, 3).Value = "হালদারসভিক গ্রামে থেকে পাওয়া যায় পরিচিত দ্বীপ এস্টুরোয়ের দৃশ্য এবং এখানে রয়েছে একটি অসাধারণ আটপাশের চর্চ।"
$ws.Cells.Item(12, 4).Value = 0.01755638228482643

$ws.Cells.Item(21, 3).Value = "সমগ্র তবেরকুলোজিস-আক্রান্ত ব্যক্তিদের মধ্যে XDR-টিবির অনুপাত কিন্তু কম মনে হচ্ছে; দক্ষিণ আফ্রিকায় কোন নির্দিষ্ট সময়ে সংক্রমিত মোট ৩,০,০০০ ব্যক্তির মধ্যে ৬,০০০ জন।"
$ws.Cells.Item(21, 4).Value = 0.03235827053239344

$ws.Cells.Item(30, 3).Value = "আমেরিকার ইরানের জাগ্রোস পর্বতমালায় প্রায় ১০,০০০ বছর আগে বাঁশগাহরুর প্রথম পালন করা হয়েছিল, তা মনে হয়।"
$ws.Cells.Item(30, 4).Value = 0.2014941615706458

$ws.Cells.Item(37, 3).Value = "টঙ্গায় কোনও বিপজ্জনক ক্ষতি বা আঘাতের কোনও খবর পর্যন্ত এখনও জানানো হয়নি, কিন্তু বিদ্যুৎ সম্পর্ক আস্তগামী ভাবে বিচ্ছিন্ন হয়ে গিয়েছিল, যার ফলে টঙ্গার কর্তৃপক্ষগুলিকে পিটিডব্লিউসি-এর প্রদত্ত সুনামী বিজ্ঞপ্তি পেতে বাধা দেওয়া হয়েছিল বলে অনুমান করা হচ্ছে।"
$ws.Cells.Item(37, 4).Value = 0.02182976214283967

$ws.Cells.Item(40, 3).Value = "নারীরা: বিবাহিত না হলেও যে কোনো মহিলা ভ্রমণকারীকে বিবাহিত বলতে বলা হয়।"
$ws.Cells.Item(40, 4).Value = 0.02233077788489958

$ws.Cells.Item(53, 3).Value = "ডাস্টিন `"গোল্ডাস্ট`" রানেল্স মন্তব্য করেছেন, `"লুনা যতটা অস্বাভাবিক ছিল আমার, সম্ভবত তার চেয়েও বেশি."
$ws.Cells.Item(53, 4).Value = 0.006209443871385971

$ws.Cells.Item(79, 3).Value = "ইজমির হ'ল তুর্কিতে তৃতীয় বৃহত্তম শহর, যেখানে বসবাসিত হ'ল প্রায় ৩৭ লাখ জন, ইস্তানবুলের পরে দ্বিতীয় বৃহত্তম বন্দর, এবং একটি খুবই উত্তম যানবাহন কেন্দ্র।"
$ws.Cells.Item(79, 4).Value = 0.1082039762100426

$ws.Cells.Item(93, 3).Value = "বাংলা অনুবাদ: ফিনল্যান্ডে জলযান চলাচল একটি জাতীয় পছন্দ, যেখানে সাত অথবা আট জনের জন্য একটি জলযান রয়েছে।"
$ws.Cells.Item(93, 4).Value = 0.04531300909655814

$ws.Cells.Item(110, 3).Value = "গ্রীষ্মকালে, খাসজীর জনপ্রিয় খাবার হচ্ছে পা অ্যাম্ব অলি: রুটি, যার সাথে দিয়ে খাওয়া হয় জলক তেল, টমেটো, এবং যে কোনও উপলব্ধ অস্ত্র যেমন চিজ, টুনাফিশ, ইত্যাদি।"
$ws.Cells.Item(110, 4).Value = 0.03189964696245855

$ws.Cells.Item(112, 3).Value = "চাকো অঞ্চলটি গুয়াইকুরু এবং পায়াগুয়ারু নামক কিছু আদিবাসী জাতিগুলির বাসস্থান ছিল, যারা শিকার করে, ফল-ফুল সংগ্রহ করে এবং মাছ পাকড়াও করে বেঁচে থাকত।"
$ws.Cells.Item(112, 4).Value = 0.11856660123276

$ws.Cells.Item(118, 3).Value = "মাইক ও'ফ্লিন, এক্টিং ডাইরেক্টর পার্ক কনসার্ভেশন অ্যান্ড হেরিটেজের মতে, প্রথম শুটিং অপারেশনের জন্য নির্বাচিত চারজন শুটারকে সংক্ষিপ্ত সুরক্ষা ও প্রশিক্ষণ নির্দেশ দেওয়া হয়েছিল।"
$ws.Cells.Item(118, 4).Value = 0.06741531344568415

$ws.Cells.Item(133, 3).Value = "একটি সংস্থা কারণে কারণে সময় সম্পদ ব্যয়কর প্রক্রিয়া ব্যবস্থাপনার জন্য অনুসরণ করতে চাইবে? সংস্থাগত শিক্ষার ধারণা ব্যাপারে অনুশীলনের একটি লক্ষ্য নবান্বেষণ।"
$ws.Cells.Item(133, 4).Value = 0.02803089189112111

$ws.Cells.Item(138, 3).Value = "বর্তমান সিনেটর এবং অর্জেন্টিনার প্রথম দাম্পত্যজীবী ক্রিস্টিনা ফার্নান্দেজ ডে কির্চনার গতকাল রাতে লা প্লাটায় তাঁর রাষ্ট্রপতি প্রার্থীত্ব ঘোষণা করেছেন, যা বুয়েনোস এয়ারেস থেকে ৫০ কিলোমিটার (৩১ মাইল) দূরে অবস্থিত।"
$ws.Cells.Item(138, 4).Value = 0.06466389207746075

$ws.Cells.Item(139, 3).Value = "তার ১,০০০তম স্টাম্প ছিল ২০০০ সালে ডেভিড ক্লফ়ার এহরেনস্ট্রাহলের মহান “গ্রেট ডিডস বাই সুইডিশ কিংস” যা গিনিস বুক অব ওয়ার্ল্ড রেকর্ডসে লিস্টেড আছে।"
$ws.Cells.Item(139, 4).Value = 0.08391878323066816

$ws.Cells.Item(141, 3).Value = "ট্রাম্প ও তুর্কি রাষ্ট্রপতি রেজেপ তাইয়িপ এর্দোগানের মধ্যে টেলিফোনে কথাবার্তার পর এ ঘোষণা করা হয়েছে।"
$ws.Cells.Item(141, 4).Value = 0.01973065769362093

$ws.Cells.Item(147, 3).Value = "তাদের সফলতার রহস্য নিহিত একটি ধারণা, যা প্রতিটি বিশেষ কাজ রাখে যা অন্য বিলাসীদের সাথে প্রতিযোগিতা থেকে বিরত থাকে।"
$ws.Cells.Item(147, 4).Value = 0.09777995826836774

$ws.Cells.Item(148, 3).Value = "সঙ্গী উস্ত্রীও লুনার প্রতি শ্রদ্ধা জ্ঞাপন করেছেন।"
$ws.Cells.Item(148, 4).Value = 0.03235934848975789
